# Horarios Linea 141 - scrape update (20:52:41 -> 22:03:17)
# Applies the scraper refresh described in the commit:
#  - sheet "LP1912": header timestamps/row-count, several tie-break
#    reorderings where the scrape time (col A) changed for identical
#    arrival times (col B), overwrite of the now-stale last batch row,
#    and 6 brand-new scraped rows appended at the end.
#  - sheet "LP1912-215": same header update, one tie-break swap, and the
#    matching new 215A_EL PATO row appended.
#  - sheet "6203-6173": header timestamp only (no new data this batch).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 22:03:17"
$ws1.Range("A3").Value = "Total filas: 374"

function Set-Row141($ws, $r, $a, $b, $c, $d, $e) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
}

# Tie-break reorders (same Hora_Llegada, scrape order changed) -------
Set-Row141 $ws1 86  "10:50:41" "10:59" "10_OLMOS"      9   "LP1912"
Set-Row141 $ws1 87  "09:23:23" "10:59" "27_EL RETIRO"  96  "LP1912"

Set-Row141 $ws1 109 "10:37:52" "11:53" "23_HERNANDEZ"  76  "LP1912"
Set-Row141 $ws1 110 "10:50:41" "11:53" "225_GOMEZ"     63  "LP1912"

$ws1.Range("C120").Value = "16_P MOR-SANTA ANA"
$ws1.Range("C121").Value = "15_ABASTO"

Set-Row141 $ws1 130 "11:47:17" "12:33" "14_ABASTO"     46  "LP1912"
Set-Row141 $ws1 131 "11:34:59" "12:33" "15_ABASTO"     59  "LP1912"

Set-Row141 $ws1 142 "11:11:33" "12:48" "15X38_ABASTO"  97  "LP1912"
Set-Row141 $ws1 143 "11:47:17" "12:48" "14_ABASTO"     61  "LP1912"

Set-Row141 $ws1 167 "12:45:56" "14:01" "23_HERNANDEZ"  76  "LP1912"
Set-Row141 $ws1 168 "12:11:52" "14:01" "10_OLMOS"      110 "LP1912"

$ws1.Range("C273").Value = "15_ABASTO"
$ws1.Range("C274").Value = "16_P MOR-SANTA ANA"

Set-Row141 $ws1 283 "17:13:39" "18:36" "23_HERNANDEZ"  83  "LP1912"
Set-Row141 $ws1 284 "16:37:06" "18:36" "15X38_ABASTO"  119 "LP1912"

Set-Row141 $ws1 297 "18:12:30" "19:10" "16_SANTA ANA"  58  "LP1912"
Set-Row141 $ws1 298 "17:56:03" "19:10" "27_EL RETIRO"  74  "LP1912"

Set-Row141 $ws1 320 "17:56:03" "19:52" "81_EL PELIGRO" 116 "LP1912"
Set-Row141 $ws1 321 "18:44:57" "19:52" "225_GOMEZ"     68  "LP1912"

# Tail of the table: the previous last batch (20:52:41 / row 372) is no
# longer the newest, so it slides down to row 373 (row 373's old content,
# which was already identical across the rename, moves to 374), the new
# 22:03:17 batch's first row takes over row 372, and its remaining rows
# are appended as brand new rows 375-379.
Set-Row141 $ws1 372 "22:03:17" "22:41" "23_HERNANDEZ"     38  "LP1912"
Set-Row141 $ws1 373 "20:52:41" "22:43" "11X44_ETCHEVERRY" 111 "LP1912"
Set-Row141 $ws1 374 "20:46:10" "22:44" "11X44_ETCHEVERRY" 118 "LP1912"
Set-Row141 $ws1 375 "22:03:17" "23:04" "15_ABASTO"        61  "LP1912"
Set-Row141 $ws1 376 "22:03:17" "23:19" "14_ABASTO"        76  "LP1912"
Set-Row141 $ws1 377 "22:03:17" "23:34" "16_SANTA ANA"     91  "LP1912"
Set-Row141 $ws1 378 "22:03:17" "23:40" "215A_EL PATO"     97  "LP1912"
Set-Row141 $ws1 379 "22:03:17" "23:59" "11X44_ETCHEVERRY" 116 "LP1912"

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 22:03:17"
$ws2.Range("A3").Value = "Total filas: 57"

$ws2.Range("C47").Value = "215B_EL PATO"
$ws2.Range("C48").Value = "215_EL PELIGRO"

Set-Row141 $ws2 62 "22:03:17" "23:40" "215A_EL PATO" 97 "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 22:03:17"
